$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- S.No / Date / Name columns for the three new rows ---
$ws.Range("B16").Value = 43497
$ws.Range("B16").NumberFormat = "mm-dd-yy"
$ws.Range("C16").Value = "D.Venkatesh"
$ws.Range("B17").Value = 43525
$ws.Range("C17").Value = "D.Venkatesh"
$ws.Range("B18").Value = 43556
$ws.Range("C18").Value = "D.Venkatesh"

# --- Assigned Project / Project Information / Remarks text, in the same
#     order the original author typed them (keeps sharedStrings ordering
#     close to the authored workbook) ---
$ws.Range("E16").Value = "Working on TTS & WhataApp system"

$ws.Range("D17").Value = "Voice Recognization system & WhatsApp Project & Twitter Data Anlaytics"
$ws.Range("D17").WrapText = $true
$ws.Range("D18").Value = "Twitter Data Anlaytics & Whats App Projoect"

$ws.Range("E17").Value = "Collecting the information of Thiruvuru for Health & Education system for Making Analytics "
$ws.Range("E17").WrapText = $true
$ws.Range("E18").Value = "Making Maps for Thiruvuru Hash tags & WhatsApp Project"
$ws.Range("E18").WrapText = $true

$ws.Range("J18").Value = "Thiruvuru Data Analytics Work has completed."
$ws.Range("J17").Value = "Thiuvuru Education & Health Information has been Completed "
$ws.Range("J17").WrapText = $true

$ws.Range("D16").Value = "Voice Recognization System & WhatsApp Project"
$ws.Range("D16").WrapText = $true

# --- office Log-in / Office Log-out time values ---
$ws.Range("H16").Value = 0.35416666666666669
$ws.Range("H16").NumberFormat = "h:mm"
$ws.Range("I16").Value = 0.20833333333333334
$ws.Range("I16").NumberFormat = "h:mm"
$ws.Range("H17").Value = 0.35416666666666669
$ws.Range("H17").NumberFormat = "h:mm"
$ws.Range("I17").Value = 0.20833333333333334
$ws.Range("I17").NumberFormat = "h:mm"
$ws.Range("H18").Value = 0.35416666666666669
$ws.Range("H18").NumberFormat = "h:mm"
$ws.Range("I18").Value = 0.20833333333333334
$ws.Range("I18").NumberFormat = "h:mm"

# --- Row heights (wrapped text rows grew taller) ---
$ws.Rows.Item(16).RowHeight = 25.5
$ws.Rows.Item(17).RowHeight = 38.25
$ws.Rows.Item(18).RowHeight = 25.5

# Re-use the same new date-number-format style for B17:B18 instead of
# letting Excel mint a second/third near-identical style entry.
$ws.Range("B16").Copy()
$ws.Range("B17:B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection moves to A17 ---
$ws.Range("A17").Select()
